# The presentation's Design theme is switched from the custom "Integral"
# theme to the built-in "Office Theme" palette (Design tab -> Themes
# gallery -> Office). This rewrites the 12 theme colour slots used by the
# slide master's theme part (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) to
# the stock "Office" values; the font scheme and format scheme (fills,
# lines, effects) are already shared between the two themes, so nothing
# else needs to change there.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = 0         # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477   # folHlink 954F72
